# Fruta / hortaliza, semanal
# Insert two new weekly records at rows 84 and 85, pushing all subsequent
# records down by two rows (old row 84 -> new row 86, ..., old row 199 -> new row 201).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 84.
$ws.Rows.Item(84).Insert()
$ws.Rows.Item(84).Insert()

# New row 84
$ws.Range("A84").Value = 4
$ws.Range("B84").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C84").Value = 'Los Lagos'
$ws.Range("D84").Value = 44579
$ws.Range("E84").Value = 10
$ws.Range("F84").Value = 100112044
$ws.Range("G84").Value = 'Perejil'
$ws.Range("H84").Value = 'Sin especificar'
$ws.Range("I84").Value = 'Primera'
$ws.Range("J84").Value = 100
$ws.Range("K84").Value = 7000
$ws.Range("L84").Value = 7000
$ws.Range("M84").Value = 7000
$ws.Range("N84").Value = '$/docena de atados (2 kilos)'
$ws.Range("O84").Value = 'Región de La Araucanía'
$ws.Range("P84").Value = 3500
$ws.Range("Q84").Value = 2
$ws.Range("R84").Value = 'Hortaliza'

# New row 85
$ws.Range("A85").Value = 4
$ws.Range("B85").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C85").Value = 'Los Lagos'
$ws.Range("D85").Value = 44579
$ws.Range("E85").Value = 10
$ws.Range("F85").Value = 100112044
$ws.Range("G85").Value = 'Perejil'
$ws.Range("H85").Value = 'Sin especificar'
$ws.Range("I85").Value = 'Primera'
$ws.Range("J85").Value = 140
$ws.Range("K85").Value = 6000
$ws.Range("L85").Value = 6000
$ws.Range("M85").Value = 6000
$ws.Range("N85").Value = '$/docena de atados (3 kilos)'
$ws.Range("O85").Value = 'Región Metropolitana'
$ws.Range("P85").Value = 2000
$ws.Range("Q85").Value = 3
$ws.Range("R85").Value = 'Hortaliza'
